$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the penjualan_kode column (previously PJ-1011 / PJ-1012)
$ws.Range("A2").Value = "PJ-1014"
$ws.Range("A3").Value = "PJ-1015"

# Update the pembeli column (previously Administrator / Administrator)
$ws.Range("C2").Value = "Djoko Susanto"
$ws.Range("C3").Value = "Bachtiar Karim"

# Update row 2 data (date + user_id)
$ws.Range("B2").Value = 45760.51394675926
$ws.Range("D2").Value = 4

# Update row 3 data (date + user_id)
$ws.Range("B3").Value = 45760.51394675926
$ws.Range("D3").Value = 4

# Remove the old row 4 (previously PJ-1013) entirely, shifting dimension to A1:G3
$ws.Rows(4).Delete()

# Update the selection to match the new state of the workbook
$ws.Range("A4:G4").Select() | Out-Null
